$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.132.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.961.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0857"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.423.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.42%  "
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.952.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.147.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("E25").Value = "  +4.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.03"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0435"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.116"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.270"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.030.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0322"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.66%  "
